# Auto-generated edit script applying the Sophia_Profits.xlsx market-price refresh diff.
# Each block updates one leve row (currentAveragePrice* / Leve*Price* / Leve*Profit* columns)
# on a specific crafting-class worksheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

# ALC row 20: Shut Up and Take My Gil
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 9673.666999999999
$ws.Range("I20").Value = 2510.5
$ws.Range("J20").Value = 24000
$ws.Range("K20").Value = 2510.5
$ws.Range("L20").Value = 24000
$ws.Range("M20").Value = -2280.5
$ws.Range("N20").Value = -24460

# ALC row 35: Conspicuous Conjuration
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H35").Value = 9673.666999999999
$ws.Range("I35").Value = 2510.5
$ws.Range("J35").Value = 24000
$ws.Range("K35").Value = 2510.5
$ws.Range("L35").Value = 24000
$ws.Range("M35").Value = -2131.5
$ws.Range("N35").Value = -24758

# ALC row 62: The Mustache Suits Him
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3533.3333
$ws.Range("I62").Value = 3533.3333
$ws.Range("K62").Value = 3533.3333
$ws.Range("M62").Value = -2909.3333

# ALC row 64: Forged from the Void
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 5999.6
$ws.Range("I64").Value = 5999.6
$ws.Range("K64").Value = 5999.6
$ws.Range("M64").Value = -5751.6

# ALC row 65: Forgery of Convenience (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 3533.3333
$ws.Range("I65").Value = 3533.3333
$ws.Range("K65").Value = 17666.6665
$ws.Range("M65").Value = -14546.6665

# ALC row 67: Dodging the Draft (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 5999.6
$ws.Range("I67").Value = 5999.6
$ws.Range("K67").Value = 5999.6
$ws.Range("M67").Value = -5141.6

# ALC row 87: There Was a Late Fee
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 81687.336
$ws.Range("J87").Value = 81687.336
$ws.Range("L87").Value = 81687.336
$ws.Range("N87").Value = -84183.336

# ALC row 90: A Gate Arcane Is Dragon's Bane (L)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 81687.336
$ws.Range("J90").Value = 81687.336
$ws.Range("L90").Value = 245062.008
$ws.Range("N90").Value = -257542.008

# ALC row 111: An Eye for Healing
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 1867.5555
$ws.Range("I111").Value = 847
$ws.Range("K111").Value = 2541
$ws.Range("M111").Value = 526

# ALC row 132: Fast-forwarding Flora
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1177.7142
$ws.Range("I132").Value = 1177.7142
$ws.Range("K132").Value = 3533.1426
$ws.Range("M132").Value = -1003.1426

# ALC row 137: Cutting Edge of Culinary Quality
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3146.16
$ws.Range("I137").Value = 3020.9412
$ws.Range("J137").Value = 3412.25
$ws.Range("K137").Value = 9062.8236
$ws.Range("L137").Value = 10236.75
$ws.Range("M137").Value = -6512.8236
$ws.Range("N137").Value = -15336.75

# ALC row 138: All-night Crafting
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1821.3125
$ws.Range("J138").Value = 2731.125
$ws.Range("L138").Value = 8193.375
$ws.Range("N138").Value = -18473.375

# ARM row 32: Ingot We Trust
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12561.107
$ws.Range("I32").Value = 12561.107
$ws.Range("K32").Value = 12561.107
$ws.Range("M32").Value = -12274.107

# ARM row 61: Dealing with the Tough Stuff
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1207.75
$ws.Range("I61").Value = 1207.75
$ws.Range("K61").Value = 1207.75
$ws.Range("M61").Value = -995.75

# ARM row 74: As the Bolt Flies
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 16662.324
$ws.Range("I74").Value = 16917.367
$ws.Range("K74").Value = 16917.367
$ws.Range("M74").Value = -16043.367

# ARM row 77: Heavy Metal Banned (L)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 16662.324
$ws.Range("I77").Value = 16917.367
$ws.Range("K77").Value = 84586.83499999999
$ws.Range("M77").Value = -80218.83499999999

# ARM row 110: Scheduled Maintenance
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 5033
$ws.Range("I110").Value = 2537.125
$ws.Range("K110").Value = 2537.125
$ws.Range("M110").Value = -492.125

# ARM row 132: Don't Bore Me, Ore Me
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2639.5557
$ws.Range("I132").Value = 2465.1428
$ws.Range("J132").Value = 3250
$ws.Range("K132").Value = 7395.428400000001
$ws.Range("L132").Value = 9750
$ws.Range("M132").Value = -4865.428400000001
$ws.Range("N132").Value = -14810

# ARM row 136: Metal with Mettle
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1207.75
$ws.Range("I136").Value = 1207.75
$ws.Range("K136").Value = 3623.25
$ws.Range("M136").Value = -1073.25

# BSM row 80: Unbreaker
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 423.30768
$ws.Range("I80").Value = 196.33333
$ws.Range("J80").Value = 491.4
$ws.Range("K80").Value = 196.33333
$ws.Range("L80").Value = 491.4
$ws.Range("M80").Value = 801.6666700000001
$ws.Range("N80").Value = -2487.4

# BSM row 83: Attack on Titanium (L)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 423.30768
$ws.Range("I83").Value = 196.33333
$ws.Range("J83").Value = 491.4
$ws.Range("K83").Value = 981.6666499999999
$ws.Range("L83").Value = 2457
$ws.Range("M83").Value = 4010.33335
$ws.Range("N83").Value = -12441

# BSM row 134: Ruthenium Supremium
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3660
$ws.Range("I134").Value = 3634.5454
$ws.Range("J134").Value = 3800
$ws.Range("K134").Value = 10903.6362
$ws.Range("L134").Value = 11400
$ws.Range("M134").Value = -8368.636200000001
$ws.Range("N134").Value = -16470

# CRP row 31: Wall Not Found
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2028.4706
$ws.Range("I31").Value = 1898.8182
$ws.Range("K31").Value = 1898.8182
$ws.Range("M31").Value = -1603.8182

# CRP row 34: Armoires of the Rich and Famous
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2028.4706
$ws.Range("I34").Value = 1898.8182
$ws.Range("K34").Value = 1898.8182
$ws.Range("M34").Value = -1696.8182

# CRP row 58: You Do the Heavy Lifting
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4943.2856
$ws.Range("I58").Value = 4785.077
$ws.Range("K58").Value = 4785.077
$ws.Range("M58").Value = -4582.077

# CRP row 70: A Reward Fitting of the Faithful
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# CRP row 73: Just Rewards for Just Devotion (L)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# CRP row 132: Hull Lotta Damage
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# CRP row 134: Wood You Be Quiet
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 13333
$ws.Range("I134").Value = 15000
$ws.Range("K134").Value = 45000
$ws.Range("M134").Value = -42465

# CRP row 136: Turali Quality
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 4943.2856
$ws.Range("I136").Value = 4785.077
$ws.Range("K136").Value = 14355.231
$ws.Range("M136").Value = -11805.231

# CUL row 87: Soup That Eats Like a Knight
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 7245.375
$ws.Range("I87").Value = 7588
$ws.Range("J87").Value = 6674.3335
$ws.Range("K87").Value = 22764
$ws.Range("L87").Value = 20023.0005
$ws.Range("M87").Value = -21516
$ws.Range("N87").Value = -22519.0005

# CUL row 90: Like Ma Used to Make (L)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H90").Value = 7245.375
$ws.Range("I90").Value = 7588
$ws.Range("J90").Value = 6674.3335
$ws.Range("K90").Value = 68292
$ws.Range("L90").Value = 60069.0015
$ws.Range("M90").Value = -62052
$ws.Range("N90").Value = -72549.0015

# GSM row 80: Needs More Prayerbell
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6642.5713
$ws.Range("I80").Value = 2833
$ws.Range("J80").Value = 9499.75
$ws.Range("K80").Value = 2833
$ws.Range("L80").Value = 9499.75
$ws.Range("M80").Value = -1835
$ws.Range("N80").Value = -11495.75

# GSM row 83: With a Noise That Reaches Heaven (L)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6642.5713
$ws.Range("I83").Value = 2833
$ws.Range("J83").Value = 9499.75
$ws.Range("K83").Value = 14165
$ws.Range("L83").Value = 47498.75
$ws.Range("M83").Value = -9173
$ws.Range("N83").Value = -57482.75

# GSM row 102: Put the Metal to the Peddle
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4697.25
$ws.Range("I102").Value = 4697.25
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 4697.25
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -3075.25
$ws.Range("N102").ClearContents()

# GSM row 132: On Board for Lar
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3600.3333
$ws.Range("I132").Value = 3520.4
$ws.Range("K132").Value = 10561.2
$ws.Range("M132").Value = -8031.200000000001

# LTW row 132: Tenets of Tanning
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4093.6667
$ws.Range("I132").Value = 3849.3333
$ws.Range("J132").Value = 4582.3335
$ws.Range("K132").Value = 11547.9999
$ws.Range("L132").Value = 13747.0005
$ws.Range("M132").Value = -9017.999899999999
$ws.Range("N132").Value = -18807.0005

# LTW row 136: Respect for Br'aax
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2987.111
$ws.Range("I136").Value = 2834.4666
$ws.Range("K136").Value = 8503.399800000001
$ws.Range("M136").Value = -5953.399800000001

# WVR row 81: Where the Dragonflies, the Net Catches
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4489.6665
$ws.Range("I81").Value = 1925.875
$ws.Range("K81").Value = 3851.75
$ws.Range("M81").Value = -2790.75

# WVR row 84: To Kill a Dragon on Nameday (L)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 4489.6665
$ws.Range("I84").Value = 1925.875
$ws.Range("K84").Value = 19258.75
$ws.Range("M84").Value = -13954.75

# WVR row 96: Skills on Display
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3306.6365
$ws.Range("I96").Value = 5098.8335
$ws.Range("J96").Value = 1156
$ws.Range("K96").Value = 5098.8335
$ws.Range("L96").Value = 1156
$ws.Range("M96").Value = -3725.8335
$ws.Range("N96").Value = -3902

# WVR row 132: Comfy Cabins
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 895.0833
$ws.Range("I132").Value = 895.0833
$ws.Range("K132").Value = 2685.2499
$ws.Range("M132").Value = -155.2498999999998

# WVR row 136: Weaving the Envelope
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1932.875
$ws.Range("I136").Value = 1950.591
$ws.Range("K136").Value = 5851.772999999999
$ws.Range("M136").Value = -3301.772999999999
